$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '73.268.13'
$ws.Cells.Item(2, 5).Value = '  +1.63%  '

$ws.Cells.Item(3, 4).Value = '3.982.69'
$ws.Cells.Item(3, 5).Value = '  -1.28%  '

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '0.999'
$ws.Cells.Item(4, 5).Value = '  -0.05%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '616.51'
$ws.Cells.Item(5, 5).Value = '  +14.31%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '167.06'
$ws.Cells.Item(6, 5).Value = '  +12.10%  '

$ws.Cells.Item(7, 5).Value = '  -1.70%  '

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '1.00'
$ws.Cells.Item(8, 5).Value = '  +0.02%  '

$ws.Cells.Item(9, 5).Value = '  +0.52%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.187'
$ws.Cells.Item(10, 5).Value = '  +8.30%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '56.59'
$ws.Cells.Item(11, 5).Value = '  +6.16%  '

$ws.Cells.Item(12, 5).Value = '  +2.48%  '

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '11.12'
$ws.Cells.Item(13, 5).Value = '  +2.52%  '

$ws.Cells.Item(14, 4).Value = '4.621.60'

$ws.Cells.Item(15, 4).Value = '3.991.96'
$ws.Cells.Item(15, 5).Value = '  -1.28%  '

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '1.25'
$ws.Cells.Item(16, 5).Value = '  +3.43%  '

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '14.13'
$ws.Cells.Item(17, 5).Value = '  -1.21%  '

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '20.52'
$ws.Cells.Item(18, 5).Value = '  -0.46%  '

$ws.Cells.Item(19, 4).Value = '73.221.12'
$ws.Cells.Item(19, 5).Value = '  +1.59%  '

$ws.Cells.Item(20, 5).Value = '  -0.31%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '440.09'
$ws.Cells.Item(21, 5).Value = '  +0.20%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '4.88'
$ws.Cells.Item(22, 5).Value = '  +14.11%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '95.73'
$ws.Cells.Item(23, 5).Value = '  -2.08%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '3.37'
$ws.Cells.Item(24, 5).Value = '  -3.61%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '14.18'
$ws.Cells.Item(25, 5).Value = '  -2.69%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '4.09'

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '11.16'
$ws.Cells.Item(27, 5).Value = '  -0.82%  '

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '5.94'
$ws.Cells.Item(28, 5).Value = '  -0.16%  '

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '10.50'
$ws.Cells.Item(29, 5).Value = '  -1.55%  '

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '36.11'
$ws.Cells.Item(30, 5).Value = '  -2.74%  '

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '7.80'
$ws.Cells.Item(31, 5).Value = '  -5.99%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '13.70'
$ws.Cells.Item(32, 5).Value = '  +1.36%  '

$ws.Cells.Item(33, 5).Value = '  -2.87%  '

$ws.Cells.Item(34, 2).Value = 'PEPE'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '0.0000101'
$ws.Cells.Item(34, 5).Value = '  +10.93%  '

$ws.Cells.Item(35, 2).Value = 'OKB'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '71.08'
$ws.Cells.Item(35, 5).Value = '  +6.66%  '

$ws.Cells.Item(36, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '47.67'
$ws.Cells.Item(36, 5).Value = '  -3.60%  '

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '638.48'
$ws.Cells.Item(37, 5).Value = '  -6.40%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.433'
$ws.Cells.Item(38, 5).Value = '  -5.45%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '3.44'
$ws.Cells.Item(39, 5).Value = '  +1.21%  '

$ws.Cells.Item(40, 2).Value = 'Kaspa'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.147'
$ws.Cells.Item(40, 5).Value = '  -0.96%  '

$ws.Cells.Item(41, 2).Value = 'Dai'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.999'
$ws.Cells.Item(41, 5).Value = '  -0.12%  '

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '10.99'
$ws.Cells.Item(42, 5).Value = '  -3.00%  '

$ws.Cells.Item(43, 5).Value = '  +0.20%  '

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '3.26'
$ws.Cells.Item(44, 5).Value = '  -3.75%  '

$ws.Cells.Item(45, 5).Value = '  -1.34%  '

$ws.Cells.Item(46, 5).Value = '  -0.67%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '3.43'
$ws.Cells.Item(47, 5).Value = '  +2.70%  '

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '2.91'
$ws.Cells.Item(48, 5).Value = '  +28.21%  '

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '2.62'
$ws.Cells.Item(49, 5).Value = '  -0.76%  '

$ws.Cells.Item(50, 4).Value = '2.851.58'
$ws.Cells.Item(50, 5).Value = '  +3.20%  '

$ws.Cells.Item(51, 2).Value = 'FLOKI'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.000280'
$ws.Cells.Item(51, 5).Value = '  -1.83%  '
